$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "001"
$ws.Range("M2").Value = "2020-12-18 00:00:00"
$ws.Range("N2").Value = "2017-12-31 00:00:00"
$ws.Range("O2").Value = 161561304.1
$ws.Range("P2").Value = 2700086378.65
$ws.Range("Q2").Value = 2528500292.11
$ws.Range("S2").Value = 1538299860.75
$ws.Range("T2").Value = 1538299860.75
$ws.Range("V2").Value = 652293211.21
$ws.Range("W2").Value = 161086114.15
$ws.Range("X2").Value = 83228822.56
$ws.Range("Y2").Value = 188948044.98
$ws.Range("Z2").Value = 199035943.41
$ws.Range("AA2").Value = 38161559.7
$ws.Range("AG2").Value = 38082779.41
$ws.Range("AS2").Value = 119631295.25
